# Select B2:E3 on Sheet2 and remove the previous tabSelected/topLeftCell state
$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("B2:E3").Select()

# Add a new Sheet3 after Sheet2 and make it the active sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Headers
$ws3.Range("B3").Value = "APP"
$ws3.Range("D3").Value = "NANO"
$ws3.Range("B4").Value = "roll"
$ws3.Range("C4").Value = "pitch"
$ws3.Range("D4").Value = "x"
$ws3.Range("E4").Value = "y"

$ws3.Range("B5").Value = 0.7
$ws3.Range("C5").Value = -0.8
$ws3.Range("D5").Value = 0.11
$ws3.Range("E5").Value = 0.18
$ws3.Range("B6").Value = 3.5
$ws3.Range("C6").Value = 2.3
$ws3.Range("D6").Value = 2.76
$ws3.Range("E6").Value = -2.8
$ws3.Range("B7").Value = 6.9
$ws3.Range("C7").Value = 5.6
$ws3.Range("D7").Value = 5.78
$ws3.Range("E7").Value = -5.9
$ws3.Range("B8").Value = 10.5
$ws3.Range("C8").Value = 9
$ws3.Range("D8").Value = 9.15
$ws3.Range("E8").Value = -9.1
$ws3.Range("B9").Value = 13.1
$ws3.Range("C9").Value = 11.5
$ws3.Range("D9").Value = 11.8
$ws3.Range("E9").Value = -11.8
$ws3.Range("B10").Value = 16.7
$ws3.Range("C10").Value = 14.5
$ws3.Range("D10").Value = 15.15
$ws3.Range("E10").Value = -14.8
$ws3.Range("B11").Value = 19.8
$ws3.Range("C11").Value = 17.6
$ws3.Range("D11").Value = 18.2
$ws3.Range("E11").Value = -17.9
$ws3.Range("B12").Value = 22.9
$ws3.Range("C12").Value = 20.6
$ws3.Range("D12").Value = 21.2
$ws3.Range("E12").Value = -20.8
$ws3.Range("B13").Value = 25.9
$ws3.Range("C13").Value = 23.8
$ws3.Range("D13").Value = 24.15
$ws3.Range("E13").Value = -24
$ws3.Range("B14").Value = 28.4
$ws3.Range("C14").Value = 26.9
$ws3.Range("D14").Value = 26.8
$ws3.Range("E14").Value = 27.1
$ws3.Range("B15").Value = 31.6
$ws3.Range("C15").Value = 29.4
$ws3.Range("D15").Value = 29.8
$ws3.Range("E15").Value = 29.8
$ws3.Range("B16").Value = 35
$ws3.Range("C16").Value = 32.8
$ws3.Range("D16").Value = 33.2
$ws3.Range("E16").Value = 33
$ws3.Range("B17").Value = 37.7
$ws3.Range("D17").Value = 35.9
$ws3.Range("B18").Value = 41.4
$ws3.Range("D18").Value = 38.9
$ws3.Range("B19").Value = 44.4
$ws3.Range("D19").Value = 41.8
$ws3.Range("B20").Value = 46.5
$ws3.Range("D20").Value = 44.8
$ws3.Range("B21").Value = 49.7
$ws3.Range("D21").Value = 47.8
$ws3.Range("B22").Value = 53.2
$ws3.Range("D22").Value = 51.2
$ws3.Range("B23").Value = 55.5
$ws3.Range("D23").Value = 53.8
$ws3.Range("B24").Value = 59
$ws3.Range("D24").Value = 57
$ws3.Range("B25").Value = 61.5
$ws3.Range("D25").Value = 59.8

$ws3.Activate()
$ws3.Range("D3:D25").Select()
